# Update gh-pages to output generated at 456a3b4
# Refresh the "想去人数" (want-to-go count, column F) figures and one
# stale cover image URL (column I) on the "展览" and "全部类型" sheets.
# The two sheets list the same events but "全部类型" has a couple of
# extra rows interleaved, so each sheet's row numbers are given
# explicitly rather than derived from a shared offset.

$wb = $excel.ActiveWorkbook

# --- "展览" sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("展览")

$ws.Cells.Item(3, 6).Value  = 505
$ws.Cells.Item(4, 6).Value  = 1480
$ws.Cells.Item(9, 6).Value  = 131
$ws.Cells.Item(10, 6).Value = 718
$ws.Cells.Item(12, 6).Value = 57
$ws.Cells.Item(13, 6).Value = 307
$ws.Cells.Item(15, 6).Value = 6302
$ws.Cells.Item(15, 9).Value = "//i0.hdslb.com/bfs/openplatform/202407/HOMmCYEq1722236602657.jpeg"
$ws.Cells.Item(18, 6).Value = 142
$ws.Cells.Item(20, 6).Value = 15120
$ws.Cells.Item(21, 6).Value = 1499
$ws.Cells.Item(22, 6).Value = 265
$ws.Cells.Item(23, 6).Value = 129
$ws.Cells.Item(25, 6).Value = 10969
$ws.Cells.Item(26, 6).Value = 726
$ws.Cells.Item(27, 6).Value = 4283
$ws.Cells.Item(28, 6).Value = 219
$ws.Cells.Item(29, 6).Value = 370
$ws.Cells.Item(30, 6).Value = 4

# --- "全部类型" sheet ---------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")

$ws.Cells.Item(3, 6).Value  = 505
$ws.Cells.Item(4, 6).Value  = 1480
$ws.Cells.Item(10, 6).Value = 131
$ws.Cells.Item(11, 6).Value = 718
$ws.Cells.Item(14, 6).Value = 57
$ws.Cells.Item(15, 6).Value = 307
$ws.Cells.Item(18, 6).Value = 6302
$ws.Cells.Item(18, 9).Value = "//i0.hdslb.com/bfs/openplatform/202407/HOMmCYEq1722236602657.jpeg"
$ws.Cells.Item(21, 6).Value = 142
$ws.Cells.Item(23, 6).Value = 15120
$ws.Cells.Item(24, 6).Value = 1499
$ws.Cells.Item(25, 6).Value = 265
$ws.Cells.Item(26, 6).Value = 129
$ws.Cells.Item(28, 6).Value = 10969
$ws.Cells.Item(29, 6).Value = 726
$ws.Cells.Item(30, 6).Value = 4283
$ws.Cells.Item(31, 6).Value = 219
$ws.Cells.Item(32, 6).Value = 370
$ws.Cells.Item(33, 6).Value = 4
